# LOQ4070.xlsx content update
# - Adds a new "Docentes responsáveis" value row (shifts rows 13-23 down to 14-24)
# - Rewrites several description cells with new/expanded Portuguese text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (old rows 13-23 become 14-24) ---
$ws.Rows.Item(13).Insert()

# The inserted row copies formatting into A13 only; clear it (target row 13 has no A cell)
$ws.Range("A13").Clear()

# Borrow the B/C formatting (wrap-text body / red-text body styles) from an existing
# two-column data row so the new B13/C13 cells pick up the correct style ids.
$ws.Range("B9:C9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 10 (Objetivos:): replace the misplaced "Docentes" name with the real objective text ---
$ws.Range("B10").Value = "Apresentar conceitos cinéticos e de fenômenos de transporte relativos a sistemas heterogêneos fluido-sólido com reações catalíticas bem como aplicações industriais de reatores heterogêneos catalíticos."
$ws.Range("C10").Value = "Apresentar conceitos cinéticos e de fenômenos de transporte relativos a sistemas heterogêneos fluido-sólido com reações catalíticas bem como aplicações industriais de reatores heterogêneos catalíticos."

# --- Row 13 (new): Docentes responsáveis value ---
$ws.Range("B13").Value = "5963230 - Leandro Gonçalves de Aguiar"
$ws.Range("C13").Value = "5963230 - Leandro Gonçalves de Aguiar"

# --- Row 14 (Programa resumido:): new short-syllabus text ---
$ws.Range("B14").Value = "1. Conceitos gerais em catálise.`n2. Tipos de sistemas catalíticos.`n3. Interação fluido-sólido.`n4. Velocidade das reações catalíticas gás-sólido.`n5. Efeitos do transporte de massa e calor externo.`n6. Transporte de massa interno.`n7. Reatores heterogêneos catalíticos.`n8. Modelos de reatores heterogêneos."
$ws.Range("C14").Value = "1. Conceitos gerais em catálise.`n2. Tipos de sistemas catalíticos.`n3. Interação fluido-sólido.`n4. Velocidade das reações catalíticas gás-sólido.`n5. Efeitos do transporte de massa e calor externo.`n6. Transporte de massa interno.`n7. Reatores heterogêneos catalíticos.`n8. Modelos de reatores heterogêneos."
# Setting the multi-line value auto-fits the row; restore the authored height (60)
$ws.Rows.Item(14).RowHeight = 60

# --- Row 16 (Programa:): new full-syllabus text ---
$ws.Range("B16").Value = "1. Conceitos gerais em catálise.`n2. Tipos de sistemas catalíticos. Propriedades dos catalisadores sólidos.`n3. Adsorção de um fluido sobre sólidos. Interação fluido-sólido.`n4. Velocidade das reações catalíticas gás-sólido.`n5. Efeitos do transporte de massa e calor externo.`n6. Transporte de massa interno.`n7. Reatores heterogêneos catalíticos.`n8. Modelos de reatores heterogêneos."
$ws.Range("C16").Value = "1. Conceitos gerais em catálise.`n2. Tipos de sistemas catalíticos. Propriedades dos catalisadores sólidos.`n3. Adsorção de um fluido sobre sólidos. Interação fluido-sólido.`n4. Velocidade das reações catalíticas gás-sólido.`n5. Efeitos do transporte de massa e calor externo.`n6. Transporte de massa interno.`n7. Reatores heterogêneos catalíticos.`n8. Modelos de reatores heterogêneos."

# --- Row 19 (Método:): description of how the course is delivered ---
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# --- Row 20 (Critério:): grading criteria ---
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# --- Row 21 (Norma de recuperação:): recovery rule ---
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# --- Row 22 (Bibliografia:): full reading list ---
$ws.Range("B22").Value = "FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3. ed.  New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York: McGraw-Hill, 1981.`nDENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F.; Bischoff, K.B. Chemical Reactor Analysis and Design. 2nd. ed. New York: John Wiley & Sons, Inc. 1990.`nTextos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas de Engenharia Química."
$ws.Range("C22").Value = "FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3. ed.  New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York: McGraw-Hill, 1981.`nDENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F.; Bischoff, K.B. Chemical Reactor Analysis and Design. 2nd. ed. New York: John Wiley & Sons, Inc. 1990.`nTextos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas de Engenharia Química."
